$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48 and 50 swap places: Algorand <-> BabyDogeCoin (plus updated price/volume)
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000115"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1149"
$ws.Range("E50").Value = "  -1.76%  "

# Price (D) and volume (E) refresh for remaining rows
$ws.Range("D2").Value = "29.325.27"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.840.48"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'240.05"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'0.6268"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07413"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "'0.2894"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "'24.73"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").Value = "'0.07732"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.834.45"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("D13").Value = "'4.975"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "'0.6761"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "'0.00001019"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "'81.92"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'6.236"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "29.331.00"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'228.76"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "'0.9997"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'158.85"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'8.459"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'0.1351"
$ws.Range("E26").Value = "  -3.51%  "
$ws.Range("D27").Value = "'17.40"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'0.06584"
$ws.Range("E28").Value = "  +15.64%  "
$ws.Range("D29").Value = "'1.447"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").Value = "'1.481"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'4.058"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'1.137"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'0.6913"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "'0.01852"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "'2.823"
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("D39").Value = "1.243.63"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").Value = "'6.743"
$ws.Range("E40").Value = "  +4.15%  "
$ws.Range("D41").Value = "'0.9334"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("D42").Value = "'0.9995"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "1.999.02"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "'100.51"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").Value = "'65.53"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "'7.038"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").Value = "'1.710"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("D49").Value = "'8.994"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D51").Value = "'0.3882"
$ws.Range("E51").Value = "  -2.24%  "
